$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FL2900"
$ws.Range("B2").Value = "S010"
$ws.Range("C2").Value = "catering"
$ws.Range("D2").Value = 5085.06
# Row 3
$ws.Range("A3").Value = "FL7875"
$ws.Range("D3").Value = 9831.780000000001
# Row 4
$ws.Range("A4").Value = "FL1344"
$ws.Range("D4").Value = 1759.58
# Row 5
$ws.Range("A5").Value = "FL5235"
$ws.Range("C5").Value = "maintenance"
$ws.Range("D5").Value = 8480.870000000001
# Row 6
$ws.Range("A6").Value = "FL2734"
$ws.Range("B6").Value = "S009"
$ws.Range("D6").Value = 11339.45
# Row 7
$ws.Range("B7").Value = "S005"
$ws.Range("C7").Value = "maintenance"
$ws.Range("D7").Value = 24988.51
# Row 8
$ws.Range("A8").Value = "FL5705"
$ws.Range("B8").Value = "S002"
$ws.Range("D8").Value = 3905.01
# Row 9
$ws.Range("A9").Value = "FL9519"
$ws.Range("B9").Value = "S009"
$ws.Range("C9").Value = "catering"
$ws.Range("D9").Value = 14628.32
# Row 10
$ws.Range("A10").Value = "FL1990"
$ws.Range("B10").Value = "S004"
$ws.Range("C10").Value = "maintenance"
$ws.Range("D10").Value = 24411.22
# Row 11
$ws.Range("A11").Value = "FL4736"
$ws.Range("B11").Value = "S006"
$ws.Range("C11").Value = "maintenance"
$ws.Range("D11").Value = 4262.9
# Row 12
$ws.Range("A12").Value = "FL3271"
$ws.Range("B12").Value = "S007"
$ws.Range("D12").Value = 20385.32
# Row 13
$ws.Range("A13").Value = "FL6349"
$ws.Range("B13").Value = "S001"
$ws.Range("C13").Value = "fuel"
$ws.Range("D13").Value = 1759.58
$ws.Range("E13").Value = 1
# Row 14
$ws.Range("A14").Value = "FL7466"
$ws.Range("B14").Value = "S001"
$ws.Range("C14").Value = "fuel"
$ws.Range("D14").Value = 1759.58
$ws.Range("E14").Value = 1
# Row 15
$ws.Range("A15").Value = "FL2792"
$ws.Range("B15").Value = "S003"
$ws.Range("C15").Value = "catering"
$ws.Range("D15").Value = 1759.58
# Row 16
$ws.Range("A16").Value = "FL2240"
$ws.Range("B16").Value = "S009"
$ws.Range("D16").Value = 4511.65
# Row 17
$ws.Range("A17").Value = "FL1541"
$ws.Range("C17").Value = "maintenance"
$ws.Range("D17").Value = 22156.97
# Row 18
$ws.Range("A18").Value = "FL5961"
$ws.Range("B18").Value = "S007"
$ws.Range("D18").Value = 13109.55
# Row 19
$ws.Range("A19").Value = "FL1832"
$ws.Range("C19").Value = "fuel"
$ws.Range("D19").Value = 1759.58
$ws.Range("E19").Value = 0
# Row 20
$ws.Range("A20").Value = "FL2948"
$ws.Range("B20").Value = "S001"
$ws.Range("C20").Value = "maintenance"
$ws.Range("D20").Value = 1759.58
# Row 21
$ws.Range("A21").Value = "FL5986"
$ws.Range("B21").Value = "S003"
$ws.Range("C21").Value = "fuel"
$ws.Range("D21").Value = 1759.58
# Row 22
$ws.Range("A22").Value = "FL5421"
$ws.Range("B22").Value = "S005"
$ws.Range("C22").Value = "maintenance"
$ws.Range("D22").Value = 7918.06
# Row 23
$ws.Range("A23").Value = "FL1832"
$ws.Range("B23").Value = "S008"
$ws.Range("C23").Value = "catering"
$ws.Range("D23").Value = 12168.11
# Row 24
$ws.Range("A24").Value = "FL5961"
$ws.Range("B24").Value = "S003"
$ws.Range("C24").Value = "fuel"
$ws.Range("D24").Value = 1759.58
# Row 25
$ws.Range("A25").Value = "FL5675"
$ws.Range("C25").Value = "maintenance"
$ws.Range("D25").Value = 23223.21
# Row 26
$ws.Range("A26").Value = "FL3883"
$ws.Range("B26").Value = "S002"
$ws.Range("C26").Value = "fuel"
$ws.Range("D26").Value = 3905.01
